# Auto-generated Excel COM-interop script to apply the Ragnarok_Profits.xlsx market-data update
# Updates numeric market-data cells (columns H-N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 72223280
$ws.Range("I19").Value = 166667410
$ws.Range("J19").Value = 25001218
$ws.Range("K19").Value = 166667410
$ws.Range("L19").Value = 25001218
$ws.Range("M19").Value = -166667235
$ws.Range("N19").Value = -25001568
$ws.Range("H28").Value = 4820.7827
$ws.Range("I28").Value = 1290.9166
$ws.Range("J28").Value = 8671.546
$ws.Range("K28").Value = 1290.9166
$ws.Range("L28").Value = 8671.546
$ws.Range("M28").Value = -805.9166
$ws.Range("N28").Value = -9641.546
$ws.Range("H70").Value = 84350670
$ws.Range("J70").Value = 1427.1428
$ws.Range("L70").Value = 4281.428400000001
$ws.Range("N70").Value = -4821.428400000001
$ws.Range("H73").Value = 84350670
$ws.Range("J73").Value = 1427.1428
$ws.Range("L73").Value = 4281.428400000001
$ws.Range("N73").Value = -6153.428400000001
$ws.Range("H74").Value = 6569.8184
$ws.Range("I74").Value = 5952.75
$ws.Range("K74").Value = 5952.75
$ws.Range("M74").Value = -5016.75
$ws.Range("H76").Value = 4213.1665
$ws.Range("I76").Value = 2997
$ws.Range("K76").Value = 2997
$ws.Range("M76").Value = -2682
$ws.Range("H77").Value = 6569.8184
$ws.Range("I77").Value = 5952.75
$ws.Range("K77").Value = 29763.75
$ws.Range("M77").Value = -25083.75
$ws.Range("H79").Value = 4213.1665
$ws.Range("I79").Value = 2997
$ws.Range("K79").Value = 2997
$ws.Range("M79").Value = -1905
$ws.Range("H86").Value = 6597.1333
$ws.Range("I86").Value = 2869
$ws.Range("K86").Value = 2869
$ws.Range("M86").Value = -1746
$ws.Range("H89").Value = 6597.1333
$ws.Range("I89").Value = 2869
$ws.Range("K89").Value = 14345
$ws.Range("M89").Value = -8729
$ws.Range("H106").Value = 8130.5
$ws.Range("I106").Value = 9957.571
$ws.Range("J106").Value = 3867.3333
$ws.Range("K106").Value = 9957.571
$ws.Range("L106").Value = 3867.3333
$ws.Range("M106").Value = -9326.571
$ws.Range("N106").Value = -5129.3333
$ws.Range("H137").Value = 1466.08
$ws.Range("I137").Value = 927.2439000000001
$ws.Range("J137").Value = 3920.7778
$ws.Range("K137").Value = 2781.7317
$ws.Range("L137").Value = 11762.3334
$ws.Range("M137").Value = -231.7317000000003
$ws.Range("N137").Value = -16862.3334
$ws.Range("H138").Value = 8390.571
$ws.Range("I138").Value = 8894
$ws.Range("J138").Value = 8253.272000000001
$ws.Range("K138").Value = 26682
$ws.Range("L138").Value = 24759.816
$ws.Range("M138").Value = -21542
$ws.Range("N138").Value = -35039.81600000001
$ws.Range("H141").Value = 15630750
$ws.Range("I141").Value = 18522354
$ws.Range("K141").Value = 55567062
$ws.Range("M141").Value = -55561882

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 336181.53
$ws.Range("I61").Value = 2186.926
$ws.Range("K61").Value = 2186.926
$ws.Range("M61").Value = -1974.926
$ws.Range("H63").Value = 4555.3335
$ws.Range("I63").Value = 4333
$ws.Range("K63").Value = 4333
$ws.Range("M63").Value = -3647
$ws.Range("H66").Value = 4555.3335
$ws.Range("I66").Value = 4333
$ws.Range("K66").Value = 21665
$ws.Range("M66").Value = -18233
$ws.Range("H74").Value = 1202.2693
$ws.Range("I74").Value = 888.7143
$ws.Range("K74").Value = 888.7143
$ws.Range("M74").Value = -14.71429999999998
$ws.Range("H77").Value = 1202.2693
$ws.Range("I77").Value = 888.7143
$ws.Range("K77").Value = 4443.5715
$ws.Range("M77").Value = -75.57150000000001
$ws.Range("H102").Value = 3211
$ws.Range("I102").Value = 2201.6365
$ws.Range("J102").Value = 5986.75
$ws.Range("K102").Value = 2201.6365
$ws.Range("L102").Value = 5986.75
$ws.Range("M102").Value = -579.6365000000001
$ws.Range("N102").Value = -9230.75
$ws.Range("H110").Value = 5300.8076
$ws.Range("I110").Value = 4984.25
$ws.Range("J110").Value = 9099.5
$ws.Range("K110").Value = 4984.25
$ws.Range("L110").Value = 9099.5
$ws.Range("M110").Value = -2939.25
$ws.Range("N110").Value = -13189.5
$ws.Range("H136").Value = 336181.53
$ws.Range("I136").Value = 2186.926
$ws.Range("K136").Value = 6560.778
$ws.Range("M136").Value = -4010.778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2253.5264
$ws.Range("I80").Value = 1739.625
$ws.Range("J80").Value = 2627.2727
$ws.Range("K80").Value = 1739.625
$ws.Range("L80").Value = 2627.2727
$ws.Range("M80").Value = -741.625
$ws.Range("N80").Value = -4623.2727
$ws.Range("H83").Value = 2253.5264
$ws.Range("I83").Value = 1739.625
$ws.Range("J83").Value = 2627.2727
$ws.Range("K83").Value = 8698.125
$ws.Range("L83").Value = 13136.3635
$ws.Range("M83").Value = -3706.125
$ws.Range("N83").Value = -23120.3635
$ws.Range("H86").Value = 3616.5715
$ws.Range("I86").Value = 1494.875
$ws.Range("J86").Value = 4922.231
$ws.Range("K86").Value = 1494.875
$ws.Range("L86").Value = 4922.231
$ws.Range("M86").Value = -371.875
$ws.Range("N86").Value = -7168.231
$ws.Range("H89").Value = 3616.5715
$ws.Range("I89").Value = 1494.875
$ws.Range("J89").Value = 4922.231
$ws.Range("K89").Value = 7474.375
$ws.Range("L89").Value = 24611.155
$ws.Range("M89").Value = -1858.375
$ws.Range("N89").Value = -35843.155
$ws.Range("H99").Value = 2977.375
$ws.Range("I99").Value = 2707.25
$ws.Range("K99").Value = 2707.25
$ws.Range("M99").Value = -1209.25
$ws.Range("H105").Value = 998742.7
$ws.Range("I105").Value = 1527039.5
$ws.Range("J105").Value = 8186.125
$ws.Range("K105").Value = 1527039.5
$ws.Range("L105").Value = 8186.125
$ws.Range("M105").Value = -1525292.5
$ws.Range("N105").Value = -11680.125
$ws.Range("H107").Value = 4418.8125
$ws.Range("I107").Value = 4857.2144
$ws.Range("J107").Value = 1350
$ws.Range("K107").Value = 4857.2144
$ws.Range("L107").Value = 1350
$ws.Range("M107").Value = -2937.2144
$ws.Range("N107").Value = -5190
$ws.Range("H134").Value = 5002550
$ws.Range("I134").Value = 2787.5
$ws.Range("J134").Value = 25001600
$ws.Range("K134").Value = 8362.5
$ws.Range("L134").Value = 75004800
$ws.Range("M134").Value = -5827.5
$ws.Range("N134").Value = -75009870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2640.7273
$ws.Range("I31").Value = 2599.5789
$ws.Range("K31").Value = 2599.5789
$ws.Range("M31").Value = -2304.5789
$ws.Range("H34").Value = 2640.7273
$ws.Range("I34").Value = 2599.5789
$ws.Range("K34").Value = 2599.5789
$ws.Range("M34").Value = -2397.5789

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 349.5
$ws.Range("I2").Value = 50.4
$ws.Range("J2").Value = 563.1429000000001
$ws.Range("K2").Value = 302.4
$ws.Range("L2").Value = 3378.8574
$ws.Range("M2").Value = -189.4
$ws.Range("N2").Value = -3604.8574
$ws.Range("H139").Value = 2614.4614
$ws.Range("I139").Value = 1357.5
$ws.Range("J139").Value = 2843
$ws.Range("K139").Value = 4072.5
$ws.Range("L139").Value = 8529
$ws.Range("M139").Value = 1067.5
$ws.Range("N139").Value = -18809
$ws.Range("H140").Value = 2798.5588
$ws.Range("I140").Value = 1158.0416
$ws.Range("K140").Value = 3474.1248
$ws.Range("M140").Value = 1705.8752

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7644.8125
$ws.Range("I70").Value = 7391
$ws.Range("J70").Value = 7797.1
$ws.Range("K70").Value = 7391
$ws.Range("L70").Value = 7797.1
$ws.Range("M70").Value = -7121
$ws.Range("N70").Value = -8337.1
$ws.Range("H73").Value = 7644.8125
$ws.Range("I73").Value = 7391
$ws.Range("J73").Value = 7797.1
$ws.Range("K73").Value = 7391
$ws.Range("L73").Value = 7797.1
$ws.Range("M73").Value = -6455
$ws.Range("N73").Value = -9669.1
$ws.Range("H80").Value = 2414.3333
$ws.Range("I80").Value = 1547
$ws.Range("J80").Value = 3628.6
$ws.Range("K80").Value = 1547
$ws.Range("L80").Value = 3628.6
$ws.Range("M80").Value = -549
$ws.Range("N80").Value = -5624.6
$ws.Range("H83").Value = 2414.3333
$ws.Range("I83").Value = 1547
$ws.Range("J83").Value = 3628.6
$ws.Range("K83").Value = 7735
$ws.Range("L83").Value = 18143
$ws.Range("M83").Value = -2743
$ws.Range("N83").Value = -28127
$ws.Range("H132").Value = 3849246.8
$ws.Range("I132").Value = 2990.5
$ws.Range("J132").Value = 12503323
$ws.Range("K132").Value = 8971.5
$ws.Range("L132").Value = 37509969
$ws.Range("M132").Value = -6441.5
$ws.Range("N132").Value = -37515029

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2419.9
$ws.Range("I61").Value = 2419.9
$ws.Range("K61").Value = 2419.9
$ws.Range("M61").Value = -2217.9
$ws.Range("H93").Value = 1636821.4
$ws.Range("I93").Value = 1506.75
$ws.Range("J93").Value = 3090434.5
$ws.Range("K93").Value = 1506.75
$ws.Range("L93").Value = 3090434.5
$ws.Range("M93").Value = -258.75
$ws.Range("N93").Value = -3092930.5
$ws.Range("H113").Value = 2419.9
$ws.Range("I113").Value = 2419.9
$ws.Range("K113").Value = 2419.9
$ws.Range("M113").Value = -249.9000000000001
$ws.Range("H132").Value = 2056.6765
$ws.Range("I132").Value = 1965.871
$ws.Range("K132").Value = 5897.613
$ws.Range("M132").Value = -3367.613
$ws.Range("H136").Value = 2378.6775
$ws.Range("I136").Value = 2301.889
$ws.Range("K136").Value = 6905.667
$ws.Range("M136").Value = -4355.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 418994.6
$ws.Range("I132").Value = 2464.1765
$ws.Range("J132").Value = 1430568.4
$ws.Range("K132").Value = 7392.529500000001
$ws.Range("L132").Value = 4291705.199999999
$ws.Range("M132").Value = -4862.529500000001
$ws.Range("N132").Value = -4296765.199999999
$ws.Range("H136").Value = 239655.16
$ws.Range("I136").Value = 7301.2285
$ws.Range("J136").Value = 1256203.6
$ws.Range("K136").Value = 21903.6855
$ws.Range("L136").Value = 3768610.8
$ws.Range("M136").Value = -19353.6855
$ws.Range("N136").Value = -3773710.8
